{"js": "// \"modified URL in word\": the paragraph labelled \"URL to GitHub\n// Repository:\" currently holds a placeholder excuse (\"I've been having a\n// lot of difficulties doing the simplest of things with GitHub. ...\").\n// Replace everything that follows the label with the actual repo URL.\nconst label = \"URL to GitHub Repository:\";\nconst newText = \" https://github.com/scollins97/week3homework\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(label) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Find the label inside its own paragraph, then grab everything in the\n  // paragraph that comes after it (this is robust to the exact wording /\n  // punctuation of the placeholder text that is being replaced).\n  const labelResults = target.search(label, { matchCase: true });\n  labelResults.load(\"text\");\n  await context.sync();\n\n  const labelRange = labelResults.items[0];\n  const afterLabel = labelRange.getRange(Word.RangeLocation.after);\n  const paragraphEnd = target.getRange(Word.RangeLocation.end);\n  const remainder = afterLabel.expandTo(paragraphEnd);\n\n  remainder.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"modified URL in word\": the paragraph labelled \"URL to GitHub\n# Repository:\" currently holds a placeholder excuse (\"I've been having a\n# lot of difficulties doing the simplest of things with GitHub. ...\").\n# Replace everything that follows the label with the actual repo URL.\n$d = $word.ActiveDocument\n\n$label = \"URL to GitHub Repository:\"\n$newText = \" https://github.com/scollins97/week3homework\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"$label*\") {\n        $target = $p.Range\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Locate just the label text inside the paragraph so we know exactly\n    # where it ends, regardless of whatever wording currently follows it.\n    $labelRange = $target.Duplicate\n    $find = $labelRange.Find\n    $find.ClearFormatting()\n    $find.Text = $label\n    $find.Execute() | Out-Null\n\n    # Range spanning from the end of the label to the end of the paragraph\n    # (i.e. everything that needs to be replaced).\n    $afterRange = $d.Range($labelRange.End, $target.End)\n    $afterRange.Text = $newText\n}\n"}
